$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '307.76'
    'E2' = '3.14%'
    'D3' = '36.23'
    'E3' = '3.61%'
    'D4' = '5.151'
    'E4' = '3.66%'
    'D5' = '0.08142'
    'E5' = '4.50%'
    'D6' = '1.928'
    'E6' = '2.60%'
    'B7' = 'GateToken'
    'C7' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D7' = '4.188'
    'E7' = '4.01%'
    'B8' = 'KuCoinToken'
    'C8' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D8' = '7.767'
    'E8' = '0.66%'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D9' = '0.9312'
    'E9' = '0.82%'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.1378'
    'E10' = '28.22%'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '0.1930'
    'E11' = '6.70%'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.09237'
    'E12' = '0.69%'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '0.03560'
    'E13' = '0.23%'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '0.09839'
    'E14' = '-0.35%'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D15' = '0.001423'
    'E15' = '1.42%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.005865'
    'E16' = '1.99%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.560'
    'E17' = '2.15%'
    'E18' = '2.13%'
    'D19' = '0.3440'
    'E19' = '-0.01%'
    'D20' = '0.1304'
    'E20' = '1.00%'
    'D21' = '4.907'
    'E21' = '-3.17%'
    'E22' = '4.06%'
    'D23' = '0.04533'
    'E23' = '-0.16%'
    'D24' = '0.001215'
    'E24' = '0.05%'
    'D25' = '0.004885'
    'E25' = '6.49%'
    'E26' = '-0.82%'
    'D39' = '0.01997'
    'E39' = '7.19%'
    'D40' = '0.04934'
    'E41' = '15.84%'
    'D42' = '0.007667'
    'E42' = '1.09%'
    'E43' = '4.70%'
    'E44' = '-0.97%'
    'D45' = '0.01063'
    'E45' = '-4.73%'
    'D46' = '0.00006520'
    'E46' = '8.41%'
    'D47' = '0.00000000750'
    'E47' = '-0.04%'
    'E48' = '0.29%'
    'E49' = '-8.67%'
    'E50' = '-0.04%'
    'E51' = '-0.04%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
